$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.665.26"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "1.847.92"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.034"
$ws.Range("E4").Value = "  +0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.45"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.030"
$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4382"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3790"
$ws.Range("E8").Value = "  +1.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07375"
$ws.Range("E9").Value = "  -0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8815"
$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "1.855.85"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.491"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.698"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.94"
$ws.Range("E16").Value = "  +2.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.036"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009044"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.029"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.44"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").Value = "27.688.12"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("E23").Value = "  +0.88%  "

$ws.Range("D24").Value = "2.079.01"
$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.050"
$ws.Range("E25").Value = "  +6.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.08"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("E27").Value = "  -0.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.993"
$ws.Range("E28").Value = "  +2.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.304"
$ws.Range("E29").Value = "  +0.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.68"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09037"
$ws.Range("E31").Value = "  -0.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7687"
$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("E34").Value = "  +4.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.543"
$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.031"
$ws.Range("E36").Value = "  +0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.150"
$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05256"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.842"
$ws.Range("E40").Value = "  +1.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5171"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1667"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.858"
$ws.Range("E43").Value = "  +3.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.721"
$ws.Range("E44").Value = "  +2.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.08"
$ws.Range("E45").Value = "  +1.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.67"
$ws.Range("E46").Value = "  +1.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06610"
$ws.Range("E47").Value = "  +4.27%  "

$ws.Range("E48").Value = "  +0.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.696"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4688"
$ws.Range("E50").Value = "  +0.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.894"
$ws.Range("E51").Value = "  -0.77%  "
